# Final Version V1 - Added the top companies in each sector with sum of funds received
#
# Fills in the previously-blank "Answer" columns across the four summary
# sheets of the Investments workbook with the computed results of the
# underlying (offline) EDA analysis.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Table -1.1 : "Understand the Data Set"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Table -1.1")
$ws1.Activate() | Out-Null

$ws1.Range("C5").Value = 90247           # unique companies in rounds2
$ws1.Range("C6").Value = 66368           # unique companies in companies file
$ws1.Range("C7").Value = "permalink"     # unique key column
$ws1.Range("C8").Value = "Y"             # companies in rounds2 not in companies?
$ws1.Range("C9").Value = 114942          # observations in master_frame

$ws1.Range("E10").Select() | Out-Null

# ---------------------------------------------------------------------
# Table - 2.1 : "Representative Values of Investments for Each of these
# Funding Types"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Table - 2.1")
$ws2.Activate() | Out-Null

$ws2.Range("C5").Value = 5000000         # representative funding - venture
$ws2.Range("C6").Value = 414906          # representative funding - angel
$ws2.Range("C7").Value = 300000          # representative funding - seed
$ws2.Range("C8").Value = 20000000        # representative funding - private equity
$ws2.Range("C9").Value = "venture"       # most suitable investment type

$ws2.Range("C7").Select() | Out-Null

# ---------------------------------------------------------------------
# Table -  3.1 : "Analysing the Top 3 English-Speaking Countries"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Table-3.1")
$ws3.Activate() | Out-Null

$ws3.Range("C5").Value = "USA"                     # top English speaking country
$ws3.Range("C6").Value = "GBR (United Kingdom)"    # second English speaking country
$ws3.Range("C7").Value = "IND (India)"             # third English speaking country

$ws3.Range("B10").Select() | Out-Null

# ---------------------------------------------------------------------
# Table - 5.1 : "Sector-wise Investment Analysis"
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Table-5.1")
$ws4.Activate() | Out-Null

# Total number of Investments (count)
$ws4.Range("C5").Value = 11945
$ws4.Range("D5").Value = 611
$ws4.Range("E5").Value = 327

# Total amount of investment (USD)
$ws4.Range("C6").Value = 106710641155
$ws4.Range("D6").Value = 5289424747
$ws4.Range("E6").Value = 2943543602

# Top Sector name (no. of investment-wise)
$ws4.Range("C7").Value = "Others"
$ws4.Range("D7").Value = "Others"
$ws4.Range("E7").Value = "Others"

# Second Sector name (no. of investment-wise)
$ws4.Range("C8").Value = "Social, Finance, Analytics, Advertising"
$ws4.Range("D8").Value = "Cleantech / Semiconductors "
$ws4.Range("E8").Value = "Social, Finance, Analytics, Advertising"

# Third Sector name (no. of investment-wise)
$ws4.Range("C9").Value = "Cleantech / Semiconductors "
$ws4.Range("D9").Value = "Social, Finance, Analytics, Advertising"
$ws4.Range("E9").Value = "News, Search and Messaging "
# (C8=Social/Finance, D8=Cleantech, E8=Social/Finance; C9=Cleantech, D9=Social/Finance, E9=News/Search)

# Number of investments in top sector (3)
$ws4.Range("C10").Value = 2923
$ws4.Range("D10").Value = 143
$ws4.Range("E10").Value = 109

# Number of investments in second sector (4)
$ws4.Range("C11").Value = 2658
$ws4.Range("D11").Value = 130
$ws4.Range("E11").Value = 60

# Number of investments in third sector (5)
$ws4.Range("C12").Value = 2347
$ws4.Range("D12").Value = 129
$ws4.Range("E12").Value = 52

# Highest-funded company in the top sector
$ws4.Range("C13").Value = "virtustream"
$ws4.Range("D13").Value = "electric-cloud "
$ws4.Range("E13").Value = "firstcry-com "

# Highest-funded company in the second best sector
$ws4.Range("C14").Value = "shotspotter "
$ws4.Range("D14").Value = "eusa-pharma"
$ws4.Range("E14").Value = "manthan-systems "

$ws4.Range("C19").Select() | Out-Null

# ---------------------------------------------------------------------
# Leave the first sheet active/selected, matching the saved workbook.
# ---------------------------------------------------------------------
$ws1.Activate() | Out-Null
